# Apply "Aggiunti alcuni grafici + consuntivo" updates to the consuntivo
# (actuals) table on Foglio1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Riga 15 - Progettista: ore Analista (colonna C) 30 -> 35
$ws.Range("C15").Value = 35

# Riga 17 - Amministratore: ore Responsabile (colonna B) 15 -> 20
$ws.Range("B17").Value = 20

# Riga 18 - Verificatore: ore Responsabile (colonna B) 45 -> 35
$ws.Range("B18").Value = 35

# Update the view selection to match the saved state (active cell B18)
$ws.Range("B18").Select()
